$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated consensus_x / consensus_y values (columns D and E)
$ws.Range("D2").Value = 6.109631711691939
$ws.Range("E2").Value = 3.180114391530244

$ws.Range("D3").Value = 1.32443318853413
$ws.Range("E3").Value = 1.128411107423878

$ws.Range("D4").Value = 7.502994929472846
$ws.Range("E4").Value = 8.301364246861349

# Recompute the derived error columns (error_x, error_y, euclidean_error)
for ($r = 2; $r -le 4; $r++) {
    $trueX = $ws.Cells.Item($r, 2).Value()
    $trueY = $ws.Cells.Item($r, 3).Value()
    $consX = $ws.Cells.Item($r, 4).Value()
    $consY = $ws.Cells.Item($r, 5).Value()

    $errX = $consX - $trueX
    $errY = $consY - $trueY
    $euclidean = [Math]::Sqrt([Math]::Pow($errX, 2) + [Math]::Pow($errY, 2))

    $ws.Cells.Item($r, 6).Value = $errX
    $ws.Cells.Item($r, 7).Value = $errY
    $ws.Cells.Item($r, 8).Value = $euclidean
}
